$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline updates. D/E columns are forced to Text format before
# assignment (then restored to the default "Normal" style) so that numeric-
# looking values such as "215.00" or "4.241" keep their original literal text
# instead of being auto-coerced into floating point numbers by Excel.

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '25.830.04'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  -1.25%  '
$ws.Cells.Item(2, 5).Style = 'Normal'

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.633.90'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  -1.20%  '
$ws.Cells.Item(3, 5).Style = 'Normal'

$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  -0.37%  '
$ws.Cells.Item(4, 5).Style = 'Normal'

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '215.00'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 5).Style = 'Normal'

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.5022'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  -1.78%  '
$ws.Cells.Item(6, 5).Style = 'Normal'

$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  -0.37%  '
$ws.Cells.Item(7, 5).Style = 'Normal'

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2565'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  -0.80%  '
$ws.Cells.Item(8, 5).Style = 'Normal'

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.06402'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  +0.00%  '
$ws.Cells.Item(9, 5).Style = 'Normal'

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '19.59'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  -1.73%  '
$ws.Cells.Item(10, 5).Style = 'Normal'

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07695'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  -1.76%  '
$ws.Cells.Item(11, 5).Style = 'Normal'

$ws.Cells.Item(12, 2).Value = 'Polkadot'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '4.241'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  -0.92%  '
$ws.Cells.Item(12, 5).Style = 'Normal'

$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.633.17'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  -1.31%  '
$ws.Cells.Item(13, 5).Style = 'Normal'

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '1.858.34'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  -1.25%  '
$ws.Cells.Item(14, 5).Style = 'Normal'

$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  -1.68%  '
$ws.Cells.Item(15, 5).Style = 'Normal'

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.0₅7911'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  -1.28%  '
$ws.Cells.Item(16, 5).Style = 'Normal'

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '63.47'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -0.64%  '
$ws.Cells.Item(17, 5).Style = 'Normal'

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '25.842.65'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  -1.29%  '
$ws.Cells.Item(18, 5).Style = 'Normal'

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '1.003'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -0.33%  '
$ws.Cells.Item(19, 5).Style = 'Normal'

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '201.90'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  -3.31%  '
$ws.Cells.Item(20, 5).Style = 'Normal'

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '4.328'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  -1.96%  '
$ws.Cells.Item(21, 5).Style = 'Normal'

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '9.925'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  -1.31%  '
$ws.Cells.Item(22, 5).Style = 'Normal'

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '5.966'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  -0.85%  '
$ws.Cells.Item(23, 5).Style = 'Normal'

$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  -0.29%  '
$ws.Cells.Item(24, 5).Style = 'Normal'

$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  +10.55%  '
$ws.Cells.Item(25, 5).Style = 'Normal'

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '141.40'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -1.15%  '
$ws.Cells.Item(26, 5).Style = 'Normal'

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.1138'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  -2.48%  '
$ws.Cells.Item(27, 5).Style = 'Normal'

$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  -0.66%  '
$ws.Cells.Item(28, 5).Style = 'Normal'

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '6.698'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  -3.89%  '
$ws.Cells.Item(29, 5).Style = 'Normal'

$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  -0.29%  '
$ws.Cells.Item(30, 5).Style = 'Normal'

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.04981'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -2.68%  '
$ws.Cells.Item(31, 5).Style = 'Normal'

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.257'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  -2.73%  '
$ws.Cells.Item(32, 5).Style = 'Normal'

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.175'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  -1.32%  '
$ws.Cells.Item(33, 5).Style = 'Normal'

$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  -1.61%  '
$ws.Cells.Item(34, 5).Style = 'Normal'

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.364'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  -0.35%  '
$ws.Cells.Item(35, 5).Style = 'Normal'

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.168.74'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  +1.32%  '
$ws.Cells.Item(36, 5).Style = 'Normal'

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.8911'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -3.95%  '
$ws.Cells.Item(37, 5).Style = 'Normal'

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.617'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  -4.74%  '
$ws.Cells.Item(38, 5).Style = 'Normal'

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.5584'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -1.69%  '
$ws.Cells.Item(39, 5).Style = 'Normal'

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.01558'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -2.00%  '
$ws.Cells.Item(40, 5).Style = 'Normal'

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '2.557'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  -0.04%  '
$ws.Cells.Item(41, 5).Style = 'Normal'

$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  -0.35%  '
$ws.Cells.Item(42, 5).Style = 'Normal'

$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  +0.65%  '
$ws.Cells.Item(43, 5).Style = 'Normal'

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.8069'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -3.24%  '
$ws.Cells.Item(44, 5).Style = 'Normal'

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '99.34'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  -0.76%  '
$ws.Cells.Item(45, 5).Style = 'Normal'

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.770.44'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  -1.19%  '
$ws.Cells.Item(46, 5).Style = 'Normal'

$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -1.66%  '
$ws.Cells.Item(47, 5).Style = 'Normal'

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.4515'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  -0.71%  '
$ws.Cells.Item(48, 5).Style = 'Normal'

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.002'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  -0.23%  '
$ws.Cells.Item(49, 5).Style = 'Normal'

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '54.66'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -1.72%  '
$ws.Cells.Item(50, 5).Style = 'Normal'

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.05081'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +0.87%  '
$ws.Cells.Item(51, 5).Style = 'Normal'
